$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.323.15"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "1.660.21"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.79"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.05"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.892.65"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "1.667.20"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.532"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.39"
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").Value = "27.303.42"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.77"
$ws.Range("E19").Value = "  +6.02%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("E22").Value = "  +8.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.45"
$ws.Range("E27").Value = "  +4.99%  "
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("E29").Value = "  +3.17%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "1.261.30"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.819"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "1.802.12"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.89"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.92"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("E51").Value = "  +0.18%  "
